# Update "想去人数" (F column) values on the two sheets that share this
# data ("展览" and "全部类型"). Rows 2-21 get updated per the commit diff.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1083
    3  = 389
    4  = 1495
    5  = 8736
    6  = 93
    8  = 645
    9  = 282
    11 = 16
    12 = 3585
    16 = 1182
    17 = 148
    18 = 1120
    19 = 307
    20 = 202
    21 = 2339
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
